# Daily attendance processing - 2025-10-12 22:18:35
# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every data row on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversedParts = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversedParts += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
